$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = [double]"0.01514828764759746"
$ws.Range("C2").Value = [double]"7.097389502863649e-05"
$ws.Range("D2").Value = [double]"0.1575252929769615"
$ws.Range("E2").Value = [double]"0.496779210170732"
$ws.Range("G2").Value = [double]"0.6695237646903196"

$ws.Range("B3").Value = [double]"3.230985683306322"
$ws.Range("C3").Value = [double]"1.667794583268128"
$ws.Range("D3").Value = [double]"0.1575252929769615"
$ws.Range("E3").Value = [double]"0.496779210170732"
$ws.Range("G3").Value = [double]"5.553084769722144"

$ws.Range("B4").Value = [double]"1.459612070389937"
$ws.Range("C4").Value = [double]"10.29869402782916"
$ws.Range("D4").Value = [double]"26.21740644021617"
$ws.Range("E4").Value = [double]"8.660232485948974"
$ws.Range("G4").Value = [double]"46.63594502438424"

$ws.Range("B5").Value = [double]"1.459612070389937"
$ws.Range("C5").Value = [double]"1.667794583268128"
$ws.Range("D5").Value = [double]"3.900430680208489"
$ws.Range("E5").Value = [double]"8.660232485948974"
$ws.Range("G5").Value = [double]"15.68806981981553"

$ws.Range("B6").Value = [double]"3.230985683306322"
$ws.Range("C6").Value = [double]"1.667794583268128"
$ws.Range("D6").Value = [double]"3.900430680208489"
$ws.Range("E6").Value = [double]"0.496779210170732"
$ws.Range("G6").Value = [double]"9.295990156953671"

$ws.Range("B7").Value = [double]"3.230985683306322"
$ws.Range("C7").Value = [double]"1.667794583268128"
$ws.Range("D7").Value = [double]"0.8054896365839992"
$ws.Range("E7").Value = [double]"8.660232485948974"
$ws.Range("G7").Value = [double]"14.36450238910742"

$ws.Range("B8").Value = [double]"3.230985683306322"
$ws.Range("C8").Value = [double]"1.667794583268128"
$ws.Range("D8").Value = [double]"0.1575252929769615"
$ws.Range("E8").Value = [double]"0.496779210170732"
$ws.Range("G8").Value = [double]"5.553084769722144"

$ws.Range("B9").Value = [double]"1.459612070389937"
$ws.Range("C9").Value = [double]"1.667794583268128"
$ws.Range("D9").Value = [double]"3.900430680208489"
$ws.Range("E9").Value = [double]"0.496779210170732"
$ws.Range("G9").Value = [double]"7.524616544037286"

